$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.958.47"
$ws.Range("E2").Value = "  +0.51%  "

# Row 3
$ws.Range("D3").Value = "1.890.56"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "0.8274"
$ws.Range("E5").Value = "  +8.16%  "

# Row 6
$ws.Range("D6").Value = "241.13"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "0.3211"
$ws.Range("E8").Value = "  +5.89%  "

# Row 9
$ws.Range("D9").Value = "26.54"
$ws.Range("E9").Value = "  +5.21%  "

# Row 10
$ws.Range("D10").Value = "0.06993"
$ws.Range("E10").Value = "  +2.62%  "

# Row 11
$ws.Range("D11").Value = "0.08028"
$ws.Range("E11").Value = "  +0.85%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.906.84"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7458"
$ws.Range("E13").Value = "  +1.84%  "

# Row 14
$ws.Range("D14").Value = "5.191"
$ws.Range("E14").Value = "  +0.74%  "

# Row 15
$ws.Range("D15").Value = "92.22"
$ws.Range("E15").Value = "  +1.52%  "

# Row 16
$ws.Range("D16").Value = "29.976.51"
$ws.Range("E16").Value = "  +0.56%  "

# Row 17
$ws.Range("D17").Value = "14.01"
$ws.Range("E17").Value = "  +2.09%  "

# Row 18
$ws.Range("D18").Value = "5.906"
$ws.Range("E18").Value = "  +0.24%  "

# Row 19
$ws.Range("D19").Value = "242.47"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20
$ws.Range("D20").Value = "0.000007752"
$ws.Range("E20").Value = "  +1.02%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.14%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.155.64"
$ws.Range("E22").Value = "  +0.39%  "

# Row 23
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").Value = "6.927"
$ws.Range("E24").Value = "  +0.47%  "

# Row 25
$ws.Range("D25").Value = "0.1574"
$ws.Range("E25").Value = "  +23.06%  "

# Row 26
$ws.Range("D26").Value = "167.95"
$ws.Range("E26").Value = "  +0.98%  "

# Row 27
$ws.Range("D27").Value = "9.160"
$ws.Range("E27").Value = "  -0.55%  "

# Row 28
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +1.19%  "

# Row 29
$ws.Range("D29").Value = "2.086"
$ws.Range("E29").Value = "  +3.20%  "

# Row 30
$ws.Range("E30").Value = "  -1.76%  "

# Row 31
$ws.Range("D31").Value = "1.513"
$ws.Range("E31").Value = "  +0.10%  "

# Row 32
$ws.Range("D32").Value = "4.244"
$ws.Range("E32").Value = "  -0.19%  "

# Row 33
$ws.Range("D33").Value = "0.05680"
$ws.Range("E33").Value = "  +8.98%  "

# Row 34
$ws.Range("D34").Value = "4.068"
$ws.Range("E34").Value = "  +0.51%  "

# Row 35
$ws.Range("D35").Value = "1.274"
$ws.Range("E35").Value = "  +2.27%  "

# Row 36
$ws.Range("D36").Value = "0.7294"
$ws.Range("E36").Value = "  +0.94%  "

# Row 37
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  +0.09%  "

# Row 38
$ws.Range("D38").Value = "0.01903"
$ws.Range("E38").Value = "  -0.22%  "

# Row 39
$ws.Range("D39").Value = "2.780"
$ws.Range("E39").Value = "  +0.32%  "

# Row 40
$ws.Range("D40").Value = "0.4392"
$ws.Range("E40").Value = "  +0.46%  "

# Row 41
$ws.Range("D41").Value = "71.73"
$ws.Range("E41").Value = "  +0.21%  "

# Row 42
$ws.Range("D42").Value = "5.940"
$ws.Range("E42").Value = "  -3.19%  "

# Row 43
$ws.Range("D43").Value = "0.8448"
$ws.Range("E43").Value = "  +1.41%  "

# Row 44
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$ws.Range("D45").Value = "1.888"
$ws.Range("E45").Value = "  +0.77%  "

# Row 46
$ws.Range("D46").Value = "100.99"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("D47").Value = "7.575"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").Value = "9.666"
$ws.Range("E48").Value = "  -0.06%  "

# Row 49
$ws.Range("D49").Value = "983.67"
$ws.Range("E49").Value = "  +8.47%  "

# Row 50
$ws.Range("D50").Value = "2.052.40"
$ws.Range("E50").Value = "  +0.55%  "

# Row 51
$ws.Range("D51").Value = "36.06"
$ws.Range("E51").Value = "  +0.10%  "
